$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the card name, mana cost, type and rules text into a single
# Python-tuple-like string in A2, then drop the now-redundant rows 3-5.
$ws.Range("A2").Value = "('Underworld Dreams', ['{B}{B}{B}', 'Enchantment', 'Whenever an opponent draws a card, Underworld Dreams deals 1 damage to that player.'])"

$ws.Range("A3:A5").EntireRow.Delete()
